$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39; this shifts the existing rows 39-121
# down to 40-122 (preserving their data and formatting), matching the
# dimension growing from A1:T121 to A1:T122.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly price record.
$ws.Cells.Item(39, 1).Value = 9
$ws.Cells.Item(39, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(39, 3).Value = "Metropolitana"
$ws.Cells.Item(39, 4).Value = 44952
$ws.Cells.Item(39, 5).Value = 13
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100101
$ws.Cells.Item(39, 8).Value = "Berries"
$ws.Cells.Item(39, 9).Value = 100101004
$ws.Cells.Item(39, 10).Value = "Frambuesa"
$ws.Cells.Item(39, 11).Value = "Sin especificar"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 260
$ws.Cells.Item(39, 14).Value = 8000
$ws.Cells.Item(39, 15).Value = 8000
$ws.Cells.Item(39, 16).Value = 8000
$ws.Cells.Item(39, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(39, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(39, 19).Value = 4000
$ws.Cells.Item(39, 20).Value = 2
